# Clarify nomenclature of f_reg and d_min (and refine f_size wording)
# in the "Hydropower plant parameters" sheet of the parameters_simulation workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hydropower plant parameters")

# --- Row 23: f_reg description ---
# New text is split into two runs: a normal run, and an italicized
# bracketed hint run.
$normalPart = "which fraction of the incoming water is allocated for regulated use "
$italicPart = "[leave empty if unsure - default determined by storage size will be used]"

$cellFreg = $ws.Range("B23")
$cellFreg.Value = $normalPart + $italicPart

$startItalic = $normalPart.Length + 1
$lenItalic = $italicPart.Length
$italicChars = $cellFreg.Characters($startItalic, $lenItalic)
$italicChars.Font.Italic = $true

# --- Row 24: d_min description ---
$ws.Range("B24").Value = "which fraction of the regulated use (line above) must be dispatched at stable level (eq. S4, S5)"

# --- Row 28: f_size description ---
$ws.Range("B28").Value = "this percentile controls the amount of allowed VRE overproduction (represents the % of time in which hydro+VRE may not exceed average ELCC)"

# Reflect the last-used selection in the sheet (cosmetic, matches the
# saved worksheet view state).
$ws.Range("B28").Select()
